$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 17.02141
$ws.Range("H2").Value = 51.06422999999999
$ws.Range("I2").Value = 0.6180409958166109
$ws.Range("J2").Value = 0.6180409958166109
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 37.09015533333334
$ws.Range("N2").Value = 111.270466
$ws.Range("O2").Value = 0.5716576577489302
$ws.Range("P2").Value = 0.5716576577489302
$ws.Range("Q2").Value = 631.3267408923533
$ws.Range("R2").Value = 5681.940668031179
$ws.Range("S2").Value = 0.3533078680613401
$ws.Range("T2").Value = 0.3533078680613401

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 17.02141
$ws.Range("H3").Value = 51.06422999999999
$ws.Range("I3").Value = 0.6180409958166109
$ws.Range("J3").Value = 0.6180409958166109
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 20.12511833333333
$ws.Range("N3").Value = 60.375355
$ws.Range("O3").Value = 0.310181445856982
$ws.Range("P3").Value = 0.310181445856982
$ws.Range("Q3").Value = 342.5578904501833
$ws.Range("R3").Value = 3083.02101405165
$ws.Range("S3").Value = 0.1917048496812853
$ws.Range("T3").Value = 0.1917048496812853

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 17.02141
$ws.Range("H4").Value = 51.06422999999999
$ws.Range("I4").Value = 0.6180409958166109
$ws.Range("J4").Value = 0.6180409958166109
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 7.666487
$ws.Range("N4").Value = 22.999461
$ws.Range("O4").Value = 0.1181608963940878
$ws.Range("P4").Value = 0.1181608963940878
$ws.Range("Q4").Value = 130.49441848667
$ws.Range("R4").Value = 1174.44976638003
$ws.Range("S4").Value = 0.0730282780739854
$ws.Range("T4").Value = 0.07302827807398539

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 8.84402
$ws.Range("H5").Value = 26.53206
$ws.Range("I5").Value = 0.3211230402077163
$ws.Range("J5").Value = 0.3211230402077163
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 37.09015533333334
$ws.Range("N5").Value = 111.270466
$ws.Range("O5").Value = 0.5716576577489302
$ws.Range("P5").Value = 0.5716576577489302
$ws.Range("Q5").Value = 328.0260755711067
$ws.Range("R5").Value = 2952.23468013996
$ws.Range("S5").Value = 0.1835724450143586
$ws.Range("T5").Value = 0.1835724450143586

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 8.84402
$ws.Range("H6").Value = 26.53206
$ws.Range("I6").Value = 0.3211230402077163
$ws.Range("J6").Value = 0.3211230402077163
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 20.12511833333333
$ws.Range("N6").Value = 60.375355
$ws.Range("O6").Value = 0.310181445856982
$ws.Range("P6").Value = 0.310181445856982
$ws.Range("Q6").Value = 177.9869490423667
$ws.Range("R6").Value = 1601.8825413813
$ws.Range("S6").Value = 0.09960640890961921
$ws.Range("T6").Value = 0.09960640890961921

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 8.84402
$ws.Range("H7").Value = 26.53206
$ws.Range("I7").Value = 0.3211230402077163
$ws.Range("J7").Value = 0.3211230402077163
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 7.666487
$ws.Range("N7").Value = 22.999461
$ws.Range("O7").Value = 0.1181608963940878
$ws.Range("P7").Value = 0.1181608963940878
$ws.Range("Q7").Value = 67.80256435774001
$ws.Range("R7").Value = 610.22307921966
$ws.Range("S7").Value = 0.03794418628373845
$ws.Range("T7").Value = 0.03794418628373845

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 1.675477666666667
$ws.Range("H8").Value = 5.026433
$ws.Range("I8").Value = 0.06083596397567289
$ws.Range("J8").Value = 0.0608359639756729
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 37.09015533333334
$ws.Range("N8").Value = 111.270466
$ws.Range("O8").Value = 0.5716576577489302
$ws.Range("P8").Value = 0.5716576577489302
$ws.Range("Q8").Value = 62.14372691419756
$ws.Range("R8").Value = 559.2935422277779
$ws.Range("S8").Value = 0.03477734467323146
$ws.Range("T8").Value = 0.03477734467323147

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 1.675477666666667
$ws.Range("H9").Value = 5.026433
$ws.Range("I9").Value = 0.06083596397567289
$ws.Range("J9").Value = 0.0608359639756729
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 20.12511833333333
$ws.Range("N9").Value = 60.375355
$ws.Range("O9").Value = 0.310181445856982
$ws.Range("P9").Value = 0.310181445856982
$ws.Range("Q9").Value = 33.71918630652389
$ws.Range("R9").Value = 303.472676758715
$ws.Range("S9").Value = 0.01887018726607749
$ws.Range("T9").Value = 0.01887018726607749

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 1.675477666666667
$ws.Range("H10").Value = 5.026433
$ws.Range("I10").Value = 0.06083596397567289
$ws.Range("J10").Value = 0.0608359639756729
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 7.666487
$ws.Range("N10").Value = 22.999461
$ws.Range("O10").Value = 0.1181608963940878
$ws.Range("P10").Value = 0.1181608963940878
$ws.Range("Q10").Value = 12.84502775029033
$ws.Range("R10").Value = 115.605249752613
$ws.Range("S10").Value = 0.007188432036363941
$ws.Range("T10").Value = 0.007188432036363941
